$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "29.518.64"
$ws.Cells.Item(2, 5).Value = "  +0.29%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.902.57"
$ws.Cells.Item(3, 5).Value = "  -0.84%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.006"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "337.12"
$ws.Cells.Item(5, 5).Value = "  +3.67%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.29%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4765"
$ws.Cells.Item(7, 5).Value = "  -1.10%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3997"
$ws.Cells.Item(8, 5).Value = "  -1.87%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.08043"
$ws.Cells.Item(9, 5).Value = "  -2.37%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.9908"
$ws.Cells.Item(10, 5).Value = "  -2.33%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "23.22"
$ws.Cells.Item(11, 5).Value = "  -0.56%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.882.84"
$ws.Cells.Item(12, 5).Value = "  -2.22%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.913"
$ws.Cells.Item(13, 5).Value = "  -2.69%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.111"
$ws.Cells.Item(14, 5).Value = "  -1.87%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -2.69%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.06816"
$ws.Cells.Item(16, 5).Value = "  -0.29%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.20%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.00001022"

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "17.33"
$ws.Cells.Item(19, 5).Value = "  -1.74%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.25%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "29.526.34"
$ws.Cells.Item(21, 5).Value = "  +0.31%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.505"
$ws.Cells.Item(22, 5).Value = "  -2.74%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "11.65"
$ws.Cells.Item(23, 5).Value = "  -0.95%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -1.03%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.115.64"
$ws.Cells.Item(25, 5).Value = "  -1.89%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "156.68"
$ws.Cells.Item(26, 5).Value = "  +0.68%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "6.527"
$ws.Cells.Item(27, 5).Value = "  -2.27%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -2.48%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.055"
$ws.Cells.Item(29, 5).Value = "  -2.83%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "119.10"
$ws.Cells.Item(30, 5).Value = "  -1.21%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.9968"
$ws.Cells.Item(31, 5).Value = "  -2.15%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.09537"
$ws.Cells.Item(32, 5).Value = "  -0.68%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.463"
$ws.Cells.Item(33, 5).Value = "  -3.76%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +0.79%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.533"
$ws.Cells.Item(35, 5).Value = "  -0.44%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.06459"
$ws.Cells.Item(36, 5).Value = "  +5.77%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.02240"
$ws.Cells.Item(37, 5).Value = "  -1.93%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.99%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.5816"
$ws.Cells.Item(39, 5).Value = "  -2.88%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "10.52"
$ws.Cells.Item(40, 5).Value = "  -3.04%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "7.724"
$ws.Cells.Item(41, 5).Value = "  -4.01%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1820"
$ws.Cells.Item(42, 5).Value = "  -1.43%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.449"
$ws.Cells.Item(43, 5).Value = "  +2.15%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.266"
$ws.Cells.Item(44, 5).Value = "  -1.26%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "12.18"
$ws.Cells.Item(45, 5).Value = "  -2.07%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.07405"
$ws.Cells.Item(46, 5).Value = "  -2.54%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5476"
$ws.Cells.Item(47, 5).Value = "  -2.19%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.953"
$ws.Cells.Item(48, 5).Value = "  -0.10%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "115.78"
$ws.Cells.Item(49, 5).Value = "  -1.97%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.378"
$ws.Cells.Item(50, 5).Value = "  -1.96%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "71.06"
$ws.Cells.Item(51, 5).Value = "  -1.69%  "

